$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKTBYY2DJ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T220041Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICW6XgVvS0%2BjmRl9v4%2FBo9UVV0kGuTTYcCEeljytxmg8AiEA8gxSn9nT4gt0AotXaK%2FCFdYuukNoxwuOOIDr%2B7nnGgUq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDFacFe1shrwbLVmQMSq6A2BEDJNDDqt6i2eQmLGtgO9aK2xD1Lp%2BGUGcwCuHmGEHk7XiagFj1AFIq88pXotd8aLFFSG9aNVPLtb8r9EQE1Y%2B2U7RHLup2FueUdkSUwjiyMYFb6%2BBDOIa%2FzM4%2Bp596JnJC7w3QsmO3H363mPm%2BfV8S9WwNg2ezi7Kkau9UjfKizWRkVVRh8Y8hfjE1AN7lPBzPucfbC97TvoeI7UFmBJxu39titePEzDWN48013riJpoi5WP7XkiJxKySQnzinhcsURhFjgVsupa9X4o%2F0s9x3n8NItbcRDAlvK9cFmuX7JVGSbbVrNLr6hpjpX6Jb5yYbqKoYaogN%2BMRjZnJbWy28W7T6W7b9w1%2FjCieB3pXju%2BIbmERk01If%2BusQZCDAkPTQ0TPkt8vY9koSj8hHXlnl5px6YZTwCeD%2B2RYjluYNhl7B8y%2BIDRuussmnlmQY8FAtgG1UhhbrE%2F0C4LJvhKNI4fjK8MDJwrlMQXpatsmb2jCbuWP1JbVcfymaMTNZExcXGzQGJKYdQ5xcS7EA8tAxahls2AE6hNyIWYH69BLvGMmGPsaIzk7Yk4laawqCmcVnjkWDxtjHVcw%2FKiNwgY6pQELfG9gjFSxTdsFMJbIeaPNQhaxTPIKSnsIGfvVhjhDHgB%2FiuuFw%2FfrnUKiz6ACZZaRFeT%2BALoQKyRR5NQ6fuOPi6dVnUcq9LGpuNvo5NIZ6FfeI8%2B8ranNC80AJNeWbKU9ezRoqHQqJRKfmYsvFGFxDXTdKluGWWpvuEE2yhC%2B1xl97R1vwv7BBuwLOHW7lJlprQ9sYyeoa6o7G2U1Hs2271hl8Pg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2006%3A00%3A41%20GMT&X-Amz-Signature=87430834f24d2c14408c6565145ef62a0f27b78cabef6ea7f5478131b76aed7f'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKTBYY2DJ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T220041Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICW6XgVvS0%2BjmRl9v4%2FBo9UVV0kGuTTYcCEeljytxmg8AiEA8gxSn9nT4gt0AotXaK%2FCFdYuukNoxwuOOIDr%2B7nnGgUq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDFacFe1shrwbLVmQMSq6A2BEDJNDDqt6i2eQmLGtgO9aK2xD1Lp%2BGUGcwCuHmGEHk7XiagFj1AFIq88pXotd8aLFFSG9aNVPLtb8r9EQE1Y%2B2U7RHLup2FueUdkSUwjiyMYFb6%2BBDOIa%2FzM4%2Bp596JnJC7w3QsmO3H363mPm%2BfV8S9WwNg2ezi7Kkau9UjfKizWRkVVRh8Y8hfjE1AN7lPBzPucfbC97TvoeI7UFmBJxu39titePEzDWN48013riJpoi5WP7XkiJxKySQnzinhcsURhFjgVsupa9X4o%2F0s9x3n8NItbcRDAlvK9cFmuX7JVGSbbVrNLr6hpjpX6Jb5yYbqKoYaogN%2BMRjZnJbWy28W7T6W7b9w1%2FjCieB3pXju%2BIbmERk01If%2BusQZCDAkPTQ0TPkt8vY9koSj8hHXlnl5px6YZTwCeD%2B2RYjluYNhl7B8y%2BIDRuussmnlmQY8FAtgG1UhhbrE%2F0C4LJvhKNI4fjK8MDJwrlMQXpatsmb2jCbuWP1JbVcfymaMTNZExcXGzQGJKYdQ5xcS7EA8tAxahls2AE6hNyIWYH69BLvGMmGPsaIzk7Yk4laawqCmcVnjkWDxtjHVcw%2FKiNwgY6pQELfG9gjFSxTdsFMJbIeaPNQhaxTPIKSnsIGfvVhjhDHgB%2FiuuFw%2FfrnUKiz6ACZZaRFeT%2BALoQKyRR5NQ6fuOPi6dVnUcq9LGpuNvo5NIZ6FfeI8%2B8ranNC80AJNeWbKU9ezRoqHQqJRKfmYsvFGFxDXTdKluGWWpvuEE2yhC%2B1xl97R1vwv7BBuwLOHW7lJlprQ9sYyeoa6o7G2U1Hs2271hl8Pg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2006%3A00%3A41%20GMT&X-Amz-Signature=d087b0c6c8e2287172a4ca74f6b988d65072e0fa75a28d2bb0cab711b34619a1'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKTBYY2DJ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T220041Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICW6XgVvS0%2BjmRl9v4%2FBo9UVV0kGuTTYcCEeljytxmg8AiEA8gxSn9nT4gt0AotXaK%2FCFdYuukNoxwuOOIDr%2B7nnGgUq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDFacFe1shrwbLVmQMSq6A2BEDJNDDqt6i2eQmLGtgO9aK2xD1Lp%2BGUGcwCuHmGEHk7XiagFj1AFIq88pXotd8aLFFSG9aNVPLtb8r9EQE1Y%2B2U7RHLup2FueUdkSUwjiyMYFb6%2BBDOIa%2FzM4%2Bp596JnJC7w3QsmO3H363mPm%2BfV8S9WwNg2ezi7Kkau9UjfKizWRkVVRh8Y8hfjE1AN7lPBzPucfbC97TvoeI7UFmBJxu39titePEzDWN48013riJpoi5WP7XkiJxKySQnzinhcsURhFjgVsupa9X4o%2F0s9x3n8NItbcRDAlvK9cFmuX7JVGSbbVrNLr6hpjpX6Jb5yYbqKoYaogN%2BMRjZnJbWy28W7T6W7b9w1%2FjCieB3pXju%2BIbmERk01If%2BusQZCDAkPTQ0TPkt8vY9koSj8hHXlnl5px6YZTwCeD%2B2RYjluYNhl7B8y%2BIDRuussmnlmQY8FAtgG1UhhbrE%2F0C4LJvhKNI4fjK8MDJwrlMQXpatsmb2jCbuWP1JbVcfymaMTNZExcXGzQGJKYdQ5xcS7EA8tAxahls2AE6hNyIWYH69BLvGMmGPsaIzk7Yk4laawqCmcVnjkWDxtjHVcw%2FKiNwgY6pQELfG9gjFSxTdsFMJbIeaPNQhaxTPIKSnsIGfvVhjhDHgB%2FiuuFw%2FfrnUKiz6ACZZaRFeT%2BALoQKyRR5NQ6fuOPi6dVnUcq9LGpuNvo5NIZ6FfeI8%2B8ranNC80AJNeWbKU9ezRoqHQqJRKfmYsvFGFxDXTdKluGWWpvuEE2yhC%2B1xl97R1vwv7BBuwLOHW7lJlprQ9sYyeoa6o7G2U1Hs2271hl8Pg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2006%3A00%3A41%20GMT&X-Amz-Signature=a3983f6e392a75765e3e39057f3f2500d48ce0268905062e6d02281033da3585'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSKTBYY2DJ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T220041Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICW6XgVvS0%2BjmRl9v4%2FBo9UVV0kGuTTYcCEeljytxmg8AiEA8gxSn9nT4gt0AotXaK%2FCFdYuukNoxwuOOIDr%2B7nnGgUq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDFacFe1shrwbLVmQMSq6A2BEDJNDDqt6i2eQmLGtgO9aK2xD1Lp%2BGUGcwCuHmGEHk7XiagFj1AFIq88pXotd8aLFFSG9aNVPLtb8r9EQE1Y%2B2U7RHLup2FueUdkSUwjiyMYFb6%2BBDOIa%2FzM4%2Bp596JnJC7w3QsmO3H363mPm%2BfV8S9WwNg2ezi7Kkau9UjfKizWRkVVRh8Y8hfjE1AN7lPBzPucfbC97TvoeI7UFmBJxu39titePEzDWN48013riJpoi5WP7XkiJxKySQnzinhcsURhFjgVsupa9X4o%2F0s9x3n8NItbcRDAlvK9cFmuX7JVGSbbVrNLr6hpjpX6Jb5yYbqKoYaogN%2BMRjZnJbWy28W7T6W7b9w1%2FjCieB3pXju%2BIbmERk01If%2BusQZCDAkPTQ0TPkt8vY9koSj8hHXlnl5px6YZTwCeD%2B2RYjluYNhl7B8y%2BIDRuussmnlmQY8FAtgG1UhhbrE%2F0C4LJvhKNI4fjK8MDJwrlMQXpatsmb2jCbuWP1JbVcfymaMTNZExcXGzQGJKYdQ5xcS7EA8tAxahls2AE6hNyIWYH69BLvGMmGPsaIzk7Yk4laawqCmcVnjkWDxtjHVcw%2FKiNwgY6pQELfG9gjFSxTdsFMJbIeaPNQhaxTPIKSnsIGfvVhjhDHgB%2FiuuFw%2FfrnUKiz6ACZZaRFeT%2BALoQKyRR5NQ6fuOPi6dVnUcq9LGpuNvo5NIZ6FfeI8%2B8ranNC80AJNeWbKU9ezRoqHQqJRKfmYsvFGFxDXTdKluGWWpvuEE2yhC%2B1xl97R1vwv7BBuwLOHW7lJlprQ9sYyeoa6o7G2U1Hs2271hl8Pg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2006%3A00%3A41%20GMT&X-Amz-Signature=04ec9a0b7a97c062181d9f67d8e69db6acc11938de601fa5f4366bdd46dac8b2'
